$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old last data row (row 9); remaining rows shift up naturally
$ws.Rows("9:9").Delete()

# Update header row (row 1): split each "base" column into mean/std pair
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# New header cells H1:L1 need the same style (bold, centered, bordered) as the
# existing header cells - copy formatting from G1
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update data rows with renamed/re-ordered algorithms and new mean/std values

# Row 2: LR
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8835358190196899
$ws.Range("D2").Value = 0.02161895173592906
$ws.Range("E2").Value = 0.8820922963780106
$ws.Range("F2").Value = 0.02800854418673326
$ws.Range("G2").Value = 0.8869493557087205
$ws.Range("H2").Value = 0.02711976540055858
$ws.Range("I2").Value = 0.891921027131783
$ws.Range("J2").Value = 0.02944385681119933
$ws.Range("K2").Value = 0.8763553906411049
$ws.Range("L2").Value = 0.02908422941733333

# Row 3: LDA
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.8925764558022623
$ws.Range("D3").Value = 0.02135460416543035
$ws.Range("E3").Value = 0.8983360911932341
$ws.Range("F3").Value = 0.03527728050291774
$ws.Range("G3").Value = 0.907471781040855
$ws.Range("H3").Value = 0.02313588899810559
$ws.Range("I3").Value = 0.914486434108527
$ws.Range("J3").Value = 0.0279788121958852
$ws.Range("K3").Value = 0.8926406926406925
$ws.Range("L3").Value = 0.03526294472009913

# Row 4: KNN
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.9346334310850439
$ws.Range("D4").Value = 0.01455722451269701
$ws.Range("E4").Value = 0.9315545136973709
$ws.Range("F4").Value = 0.02321273793292695
$ws.Range("G4").Value = 0.9392318449705324
$ws.Range("H4").Value = 0.02150374760195192
$ws.Range("I4").Value = 0.9440528100775193
$ws.Range("J4").Value = 0.01753936613608477
$ws.Range("K4").Value = 0.9179447536590393
$ws.Range("L4").Value = 0.03601592618840947

# Row 5: DTREE
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7974025974025973
$ws.Range("D5").Value = 0.02706873865577115
$ws.Range("E5").Value = 0.7778176135318992
$ws.Range("F5").Value = 0.04381484333519634
$ws.Range("G5").Value = 0.7696334032564179
$ws.Range("H5").Value = 0.02914658848153883
$ws.Range("I5").Value = 0.7698885658914728
$ws.Range("J5").Value = 0.03423968978427551
$ws.Range("K5").Value = 0.7903009688723974
$ws.Range("L5").Value = 0.06971292780013427

# Row 6: RTREE
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8744407205697527
$ws.Range("D6").Value = 0.02234873208954248
$ws.Range("E6").Value = 0.8746690568119139
$ws.Range("F6").Value = 0.03070906582695579
$ws.Range("G6").Value = 0.8537458795325141
$ws.Range("H6").Value = 0.02485249394820985
$ws.Range("I6").Value = 0.8701671511627908
$ws.Range("J6").Value = 0.0406830972705518
$ws.Range("K6").Value = 0.8480416408987838
$ws.Range("L6").Value = 0.03664838865895326

# Row 7: XTREE
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8336908253037286
$ws.Range("D7").Value = 0.02076211676289512
$ws.Range("E7").Value = 0.8116565545136973
$ws.Range("F7").Value = 0.03535676183866732
$ws.Range("G7").Value = 0.8056937368894216
$ws.Range("H7").Value = 0.02908184860042634
$ws.Range("I7").Value = 0.8165394864341085
$ws.Range("J7").Value = 0.05245009741893358
$ws.Range("K7").Value = 0.824757781900639
$ws.Range("L7").Value = 0.04984814205342435

# Row 8: SVM
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8906116464180979
$ws.Range("D8").Value = 0.01570398931711488
$ws.Range("E8").Value = 0.8990255561684133
$ws.Range("F8").Value = 0.03409667336054903
$ws.Range("G8").Value = 0.897547697532714
$ws.Range("H8").Value = 0.02341232496872762
$ws.Range("I8").Value = 0.9090419089147286
$ws.Range("J8").Value = 0.03009448915103605
$ws.Range("K8").Value = 0.884508348794063
$ws.Range("L8").Value = 0.02833515891435427
